$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 218-222 (columns B, E, G - and F for rows 220/221) ---

# Row 218 (Serie index 222 -> "01-01-2021")
$ws.Range("B218").Value = 2097
$ws.Range("E218").Value = 2026
$ws.Range("G218").Value = 1297

# Row 219 (Serie index 223 -> "01-02-2021")
$ws.Range("B219").Value = 2122
$ws.Range("E219").Value = 2041
$ws.Range("G219").Value = 1292

# Row 220 (Serie index 224 -> "01-03-2021")
$ws.Range("B220").Value = 2125
$ws.Range("E220").Value = 2041
$ws.Range("F220").Value = 742
$ws.Range("G220").Value = 1299

# Row 221 (Serie index 225 -> "01-04-2021")
$ws.Range("B221").Value = 2131
$ws.Range("E221").Value = 2047
$ws.Range("F221").Value = 751
$ws.Range("G221").Value = 1296

# Row 222 (Serie index 226 -> "01-05-2021")
$ws.Range("B222").Value = 2100
$ws.Range("E222").Value = 2011
$ws.Range("G222").Value = 1281

# --- Append new row 223 for "01-06-2021" ---
# Entering "01-06-2021" directly via .Value would be auto-recognized as a
# date and stored as a numeric serial (with a new number-format style).
# Instead, enter it as a text formula (="01-06-2021") so it evaluates to a
# plain string, then Copy/PasteSpecial values-only to collapse it down to a
# static text cell - this keeps it a shared string with the default
# (unstyled) cell format, matching the rest of column A.
$ws.Range("A223").Formula = "=""01-06-2021"""
$ws.Range("A223").Copy()
$ws.Range("A223").PasteSpecial(-4163)

$ws.Range("B223").Value = 2081
$ws.Range("C223").Value = 107
$ws.Range("D223").Value = 107
$ws.Range("E223").Value = 1974
$ws.Range("F223").Value = 683
$ws.Range("G223").Value = 1291
